$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I40").Value = 2375
$ws.Range("J40").Value = 2133.3333
$ws.Range("K40").Value = 2375
$ws.Range("L40").Value = 2133.3333
$ws.Range("M40").Value = -2200
$ws.Range("N40").Value = -2483.3333
$ws.Range("H100").Value = 1908.625
$ws.Range("I100").Value = 1752.7142
$ws.Range("K100").Value = 1752.7142
$ws.Range("M100").Value = -1211.7142
$ws.Range("H103").Value = 1188.8889
$ws.Range("I103").Value = 550
$ws.Range("J103").Value = 1371.4286
$ws.Range("K103").Value = 1650
$ws.Range("L103").Value = 4114.2858
$ws.Range("M103").Value = -1064
$ws.Range("N103").Value = -5286.2858
$ws.Range("H107").Value = 656.1667
$ws.Range("I107").Value = 398.33334
$ws.Range("K107").Value = 398.33334
$ws.Range("M107").Value = 1521.66666
$ws.Range("H137").Value = 30683.234
$ws.Range("I137").Value = 1226.9259
$ws.Range("J137").Value = 144300.42
$ws.Range("K137").Value = 3680.7777
$ws.Range("L137").Value = 432901.26
$ws.Range("M137").Value = -1130.7777
$ws.Range("N137").Value = -438001.26

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 11628406
$ws.Range("J2").Value = 999
$ws.Range("L2").Value = 999
$ws.Range("N2").Value = -1225
$ws.Range("H5").Value = 196.66667
$ws.Range("J5").Value = 70
$ws.Range("L5").Value = 70
$ws.Range("N5").Value = -294
$ws.Range("H26").Value = 24007
$ws.Range("I26").Value = 24007
$ws.Range("K26").Value = 24007
$ws.Range("M26").Value = -23677
$ws.Range("H32").Value = 3440.985
$ws.Range("I32").Value = 2608.6724
$ws.Range("K32").Value = 2608.6724
$ws.Range("M32").Value = -2321.6724
$ws.Range("H45").Value = 1659.2
$ws.Range("I45").Value = 1099.6666
$ws.Range("K45").Value = 1099.6666
$ws.Range("M45").Value = -722.6666
$ws.Range("H102").Value = 1780.3572
$ws.Range("I102").Value = 1493.8334
$ws.Range("J102").Value = 3499.5
$ws.Range("K102").Value = 1493.8334
$ws.Range("L102").Value = 3499.5
$ws.Range("M102").Value = 128.1666
$ws.Range("N102").Value = -6743.5
$ws.Range("H116").Value = 11628406
$ws.Range("J116").Value = 999
$ws.Range("L116").Value = 999
$ws.Range("N116").Value = -5587
$ws.Range("H122").Value = 1999.3334
$ws.Range("I122").Value = 1999.3334
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5998.0002
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3548.0002
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 2913.1333
$ws.Range("I132").Value = 3277.5
$ws.Range("J132").Value = 2780.6365
$ws.Range("K132").Value = 9832.5
$ws.Range("L132").Value = 8341.9095
$ws.Range("M132").Value = -7302.5
$ws.Range("N132").Value = -13401.9095
$ws.Range("H135").Value = 18500
$ws.Range("J135").Value = 18500
$ws.Range("L135").Value = 18500
$ws.Range("N135").Value = -28640
$ws.Range("H139").Value = 49000
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 49000
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 49000
$ws.Range("M139").ClearContents()
$ws.Range("N139").Value = -59280

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 11628406
$ws.Range("J3").Value = 999
$ws.Range("L3").Value = 999
$ws.Range("N3").Value = -1227
$ws.Range("H4").Value = 196.66667
$ws.Range("J4").Value = 70
$ws.Range("L4").Value = 70
$ws.Range("N4").Value = -300
$ws.Range("H20").Value = 2146.5
$ws.Range("I20").Value = 1957.8462
$ws.Range("K20").Value = 1957.8462
$ws.Range("M20").Value = -1710.8462
$ws.Range("H94").Value = 424.4516
$ws.Range("I94").Value = 324.89285
$ws.Range("K94").Value = 324.89285
$ws.Range("M94").Value = 126.10715
$ws.Range("H107").Value = 1404.75
$ws.Range("I107").Value = 1404.75
$ws.Range("K107").Value = 1404.75
$ws.Range("M107").Value = 515.25

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1983.25
$ws.Range("I31").Value = 1330.25
$ws.Range("J31").Value = 2636.25
$ws.Range("K31").Value = 1330.25
$ws.Range("L31").Value = 2636.25
$ws.Range("M31").Value = -1035.25
$ws.Range("N31").Value = -3226.25
$ws.Range("H34").Value = 1983.25
$ws.Range("I34").Value = 1330.25
$ws.Range("J34").Value = 2636.25
$ws.Range("K34").Value = 1330.25
$ws.Range("L34").Value = 2636.25
$ws.Range("M34").Value = -1128.25
$ws.Range("N34").Value = -3040.25
$ws.Range("H86").Value = 2211.875
$ws.Range("I86").Value = 1622.25
$ws.Range("J86").Value = 2801.5
$ws.Range("K86").Value = 1622.25
$ws.Range("L86").Value = 2801.5
$ws.Range("M86").Value = -499.25
$ws.Range("N86").Value = -5047.5
$ws.Range("H89").Value = 2211.875
$ws.Range("I89").Value = 1622.25
$ws.Range("J89").Value = 2801.5
$ws.Range("K89").Value = 8111.25
$ws.Range("L89").Value = 14007.5
$ws.Range("M89").Value = -2495.25
$ws.Range("N89").Value = -25239.5
$ws.Range("H99").Value = 1002144
$ws.Range("I99").Value = 2500756.5
$ws.Range("K99").Value = 2500756.5
$ws.Range("M99").Value = -2499258.5
$ws.Range("H126").Value = 1002144
$ws.Range("I126").Value = 2500756.5
$ws.Range("K126").Value = 7502269.5
$ws.Range("M126").Value = -7499799.5
$ws.Range("H132").Value = 3174
$ws.Range("I132").Value = 1562.1666
$ws.Range("J132").Value = 4141.1
$ws.Range("K132").Value = 4686.4998
$ws.Range("L132").Value = 12423.3
$ws.Range("M132").Value = -2156.4998
$ws.Range("N132").Value = -17483.3

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 14919.857
$ws.Range("I34").Value = 20090
$ws.Range("J34").Value = 1994.5
$ws.Range("K34").Value = 60270
$ws.Range("L34").Value = 5983.5
$ws.Range("M34").Value = -60186
$ws.Range("N34").Value = -6151.5
$ws.Range("H48").Value = 2066
$ws.Range("J48").Value = 2066
$ws.Range("L48").Value = 6198
$ws.Range("N48").Value = -6698
$ws.Range("H87").Value = 10745.25
$ws.Range("I87").Value = 5993.6665
$ws.Range("K87").Value = 17980.9995
$ws.Range("M87").Value = -16732.9995
$ws.Range("H90").Value = 10745.25
$ws.Range("I90").Value = 5993.6665
$ws.Range("K90").Value = 53942.9985
$ws.Range("M90").Value = -47702.9985
$ws.Range("H113").Value = 167582.17
$ws.Range("I113").Value = 1000003
$ws.Range("J113").Value = 1098
$ws.Range("K113").Value = 3000009
$ws.Range("L113").Value = 3294
$ws.Range("M113").Value = -2997839
$ws.Range("N113").Value = -7634
$ws.Range("H131").Value = 839.64
$ws.Range("J131").Value = 846.0205999999999
$ws.Range("L131").Value = 2538.0618
$ws.Range("N131").Value = -12618.0618

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2172.4075
$ws.Range("I102").Value = 2145.476
$ws.Range("K102").Value = 2145.476
$ws.Range("M102").Value = -523.4760000000001
$ws.Range("H113").Value = 1074
$ws.Range("J113").Value = 1499
$ws.Range("L113").Value = 1499
$ws.Range("N113").Value = -5839

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3684.75
$ws.Range("I7").Value = 2782.5715
$ws.Range("K7").Value = 2782.5715
$ws.Range("M7").Value = -2670.5715
$ws.Range("H16").Value = 2845.7778
$ws.Range("I16").Value = 2663.2307
$ws.Range("J16").Value = 3320.4
$ws.Range("K16").Value = 2663.2307
$ws.Range("L16").Value = 3320.4
$ws.Range("M16").Value = -2493.2307
$ws.Range("N16").Value = -3660.4
$ws.Range("H61").Value = 1445.2307
$ws.Range("I61").Value = 1179
$ws.Range("K61").Value = 1179
$ws.Range("M61").Value = -977
$ws.Range("H93").Value = 566.4
$ws.Range("I93").Value = 582
$ws.Range("J93").Value = 530
$ws.Range("K93").Value = 582
$ws.Range("L93").Value = 530
$ws.Range("M93").Value = 666
$ws.Range("N93").Value = -3026
$ws.Range("H100").Value = 1200
$ws.Range("I100").Value = 1650
$ws.Range("J100").Value = 300
$ws.Range("K100").Value = 1650
$ws.Range("L100").Value = 300
$ws.Range("M100").Value = -1109
$ws.Range("N100").Value = -1382
$ws.Range("H113").Value = 1445.2307
$ws.Range("I113").Value = 1179
$ws.Range("K113").Value = 1179
$ws.Range("M113").Value = 991
$ws.Range("H126").Value = 3684.75
$ws.Range("I126").Value = 2782.5715
$ws.Range("K126").Value = 8347.7145
$ws.Range("M126").Value = -5877.7145

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 3380003.2
$ws.Range("H96").Value = 2600
$ws.Range("J96").Value = 2600
$ws.Range("L96").Value = 2600
$ws.Range("N96").Value = -5346
$ws.Range("H100").Value = 708.25
$ws.Range("I100").Value = 449
$ws.Range("J100").Value = 967.5
$ws.Range("K100").Value = 898
$ws.Range("L100").Value = 1935
$ws.Range("M100").Value = -357
$ws.Range("N100").Value = -3017
$ws.Range("H113").Value = 517.8570999999999
$ws.Range("I113").Value = 411.53845
$ws.Range("J113").Value = 1900
$ws.Range("K113").Value = 1234.61535
$ws.Range("L113").Value = 5700
$ws.Range("M113").Value = 935.38465
$ws.Range("N113").Value = -10040
$ws.Range("H126").Value = 13255.637
$ws.Range("J126").Value = 10000
$ws.Range("L126").Value = 30000
$ws.Range("N126").Value = -34940
